$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header capitalization: "Big Stockroom" -> "Big StockRoom"
$ws.Range("C1").Value = "Big StockRoom"

# Add new "Cooler" header column (copy formatting from the neighboring header cell)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Cooler"

# Update row 2 (Oil): Quantity 2 -> 1, add empty Cooler cell
$ws.Range("B2").Value = 1
$ws.Range("E2").Value = ""

# Delete rows 4, 5, 6 entirely (remove the extra Oil/ketchup rows), shifting cells up
$ws.Rows("4:6").Delete()

# Update row 3: rename ketchup -> Hashbrowns, clear Stockroom Y mark, set Cooler Y mark
$ws.Range("A3").Value = "Hashbrowns"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Y"
